$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark completion checkmarks for the "SnipHub" project column (E) on the
# migration-related competence rows, mirroring the existing checkmarks
# already present in columns C/D for those rows.
$ws.Range("E27").Value = "✓"
$ws.Range("E29").Value = "✓"
$ws.Range("E31").Value = "✓"

# Update the active selection to reflect where the author left off editing.
$ws.Range("E7").Select()
